$d = $word.ActiveDocument

# 1. Merge runs in first paragraph (intro text) - no visible text change needed,
#    but ensure content matches (Find/Replace won't change run structure by itself
#    in a way that matters since text content is identical). We instead directly
#    set paragraph range text to force a single run.
$d.Content.Find.Execute("The Student Placement Eligibility Application helps students and academic institutions assess student readiness for campus placements using academic, technical, and soft skill data. The application provides personalized eligibility scores and insights via SQL queries and visualizations.", $true, $false, $false, $false, $false, $true, 1, $false, "The Student Placement Eligibility Application helps students and academic institutions assess student readiness for campus placements using academic, technical, and soft skill data. The application provides personalized eligibility scores and insights via SQL queries and visualizations.", 2) | Out-Null

# 2. Merge "Each table is exported to CSV and optionally" + " inserted into "
$d.Content.Find.Execute("Each table is exported to CSV and optionally inserted into ", $true, $false, $false, $false, $false, $true, 1, $false, "Each table is exported to CSV and optionally inserted into ", 2) | Out-Null

# 3. Merge "Sideba" + "r navigation" -> "Sidebar navigation"
$d.Content.Find.Execute("Sidebar navigation", $true, $false, $false, $false, $false, $true, 1, $false, "Sidebar navigation", 2) | Out-Null

# 4. "Gender-wise Placement Count" -> "Gender-wise Placement"
$d.Content.Find.Execute("Gender-wise Placement Count", $true, $false, $false, $false, $false, $true, 1, $false, "Gender-wise Placement", 2) | Out-Null

# 5. "Average Age by Course" -> "Average Age by Batch"
$d.Content.Find.Execute("Average Age by Course", $true, $false, $false, $false, $false, $true, 1, $false, "Average Age by Batch", 2) | Out-Null

# 6. "Students per City" -> "Top Cities by Students"
$d.Content.Find.Execute("Students per City", $true, $false, $false, $false, $false, $true, 1, $false, "Top Cities by Students", 2) | Out-Null

# 7. "Top 5 Students by Programming Problems Solved" -> "Top Students by Problems Solved"
$d.Content.Find.Execute("Top 5 Students by Programming Problems Solved", $true, $false, $false, $false, $false, $true, 1, $false, "Top Students by Problems Solved", 2) | Out-Null

# 8. "Average Soft Skills Score by Batch" -> "Average Soft Skills by Batch"
$d.Content.Find.Execute("Average Soft Skills Score by Batch", $true, $false, $false, $false, $false, $true, 1, $false, "Average Soft Skills by Batch", 2) | Out-Null

# 9. "Students with Maximum Internships Completed" -> "Top Internships Completed"
$d.Content.Find.Execute("Students with Maximum Internships Completed", $true, $false, $false, $false, $false, $true, 1, $false, "Top Internships Completed", 2) | Out-Null

# 10. "Placement Package Distribution" -> "Package Distribution"
$d.Content.Find.Execute("Placement Package Distribution", $true, $false, $false, $false, $false, $true, 1, $false, "Package Distribution", 2) | Out-Null

# 11. "Students Not Placed by Batch" -> "Not Placed by Batch"
$d.Content.Find.Execute("Students Not Placed by Batch", $true, $false, $false, $false, $false, $true, 1, $false, "Not Placed by Batch", 2) | Out-Null

# 12. "Users input academic percentages, certifications, programming stats, and soft skills." -> "Users input academic percentages, certifications, programming stats, and soft skills etc."
$d.Content.Find.Execute("Users input academic percentages, certifications, programming stats, and soft skills.", $true, $false, $false, $false, $false, $true, 1, $false, "Users input academic percentages, certifications, programming stats, and soft skills etc.", 2) | Out-Null

# 13. "Displays score with feedback (High/Moderate/Low placement chance)." -> "Displays score with feedback (Eligible/Moderate/Not Eligible)."
$d.Content.Find.Execute("Displays score with feedback (High/Moderate/Low placement chance).", $true, $false, $false, $false, $false, $true, 1, $false, "Displays score with feedback (Eligible/Moderate/Not Eligible).", 2) | Out-Null

# 14. Merge "Displays developer information a" + "nd project details."
$d.Content.Find.Execute("Displays developer information and project details.", $true, $false, $false, $false, $false, $true, 1, $false, "Displays developer information and project details.", 2) | Out-Null

# 15. "Placement trends, average packages, gender & batch analysis." -> "Placement, average packages, gender & batch analysis."
$d.Content.Find.Execute("Placement trends, average packages, gender & batch analysis.", $true, $false, $false, $false, $false, $true, 1, $false, "Placement, average packages, gender & batch analysis.", 2) | Out-Null

# 16. Merge "Placement distribution varies across batches and co" + "mpanies."
$d.Content.Find.Execute("Placement distribution varies across batches and companies.", $true, $false, $false, $false, $false, $true, 1, $false, "Placement distribution varies across batches and companies.", 2) | Out-Null

# 17. Merge "Dep" + "loy automated dashboards for universities."
$d.Content.Find.Execute("Deploy automated dashboards for universities.", $true, $false, $false, $false, $false, $true, 1, $false, "Deploy automated dashboards for universities.", 2) | Out-Null
